$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The sheet originally held 8 data rows (rows 2-9), each tagged with the
# indicator "Dias com manutenção da equipe mínima da Central de Regulação do
# SAMU". The edit splits every original row into two: one row keeps the
# (renamed) indicator "Taxa de Adesao Oficial" and a new sibling row carries
# "Taxa de Adesao" - both rows share the same Codigo IBGE / Beneficiario /
# CNPJ / value columns as the original row did.
# ---------------------------------------------------------------------------

# Capture the original 8 data rows (B:E) before we start overwriting cells.
$rowCount = 8
$orig = @()
for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $i + 2
    $orig += ,@($ws.Cells.Item($r, 2).Text, $ws.Cells.Item($r, 3).Text, $ws.Cells.Item($r, 4).Text, $ws.Cells.Item($r, 5).Value2)
}

# Write the 16 expanded rows (2 per original row) starting at row 2.
# Column B (IBGE code) is wrapped in a one-element array when assigned so it
# is stored as text rather than auto-converted to a number (it keeps the
# original cell General format / shared-string type instead of gaining a
# quote-prefix style).
$target = 2
for ($i = 0; $i -lt $rowCount; $i++) {
    $vals = $orig[$i]

    $ws.Cells.Item($target, 1).Value = "Taxa de Adesao Oficial"
    $ws.Cells.Item($target, 2).Value = ,@($vals[0])
    $ws.Cells.Item($target, 3).Value = $vals[1]
    $ws.Cells.Item($target, 4).Value = $vals[2]
    $ws.Cells.Item($target, 5).Value = $vals[3]
    $target = $target + 1

    $ws.Cells.Item($target, 1).Value = "Taxa de Adesao"
    $ws.Cells.Item($target, 2).Value = ,@($vals[0])
    $ws.Cells.Item($target, 3).Value = $vals[1]
    $ws.Cells.Item($target, 4).Value = $vals[2]
    $ws.Cells.Item($target, 5).Value = $vals[3]
    $target = $target + 1
}

# Column A best-fits the new (shorter) indicator text.
$ws.Columns.Item(1).ColumnWidth = 21.140625

# Move the active selection off the data, matching the saved view state.
$ws.Range("F2").Select()
